# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" and "全部类型" worksheets, which contain duplicated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    2 = 350
    3 = 93
    4 = 1532
    6 = 50
    7 = 130
    8 = 54
    9 = 349
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
